# Apply the commit's changes:
#  - "se arreglo la equivalencia de ramos con sus pesos correcto y de los
#    electivos (revisado con 2018)" -> the Equivalencias lookup table
#    (sheet "Equivalencias", range A1:B11) is cleared out (values removed,
#    formatting kept) now that the equivalence data lives/derives elsewhere.
#  - The Equivalencias sheet becomes the active/selected sheet in the
#    workbook (tab selection moves from MallaCurricular2020 to Equivalencias).

$wb = $excel.ActiveWorkbook

$wsEquiv = $wb.Worksheets.Item("Equivalencias")

# Clear out the stale equivalence table contents (keeps cell styles/borders).
$wsEquiv.Range("A1:B11").ClearContents()

# Make Equivalencias the active/selected sheet (matches the saved workbook
# view: tabSelected moves off MallaCurricular2020 onto Equivalencias).
$wsEquiv.Activate()
